$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 130.25
$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 21
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 91
$ws.Range("N6").Value = -1724
$ws.Range("H17").Value = 1009.625
$ws.Range("J17").Value = 1009.625
$ws.Range("L17").Value = 3028.875
$ws.Range("N17").Value = -3364.875
$ws.Range("H31").Value = 1112.6
$ws.Range("I31").Value = 1112.6
$ws.Range("K31").Value = 3337.8
$ws.Range("M31").Value = -3107.8
$ws.Range("H58").Value = 396.66666
$ws.Range("I58").Value = 400
$ws.Range("J58").Value = 393.33334
$ws.Range("K58").Value = 1200
$ws.Range("L58").Value = 1180.00002
$ws.Range("M58").Value = -1050
$ws.Range("N58").Value = -1480.00002
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = $null
$ws.Range("H138").Value = 3970
$ws.Range("J138").Value = 3953.7273
$ws.Range("L138").Value = 11861.1819
$ws.Range("N138").Value = -22141.1819
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 694.4167
$ws.Range("I2").Value = 622.2
$ws.Range("K2").Value = 622.2
$ws.Range("M2").Value = -509.2
$ws.Range("H32").Value = 2266919.2
$ws.Range("I32").Value = 2123492.8
$ws.Range("K32").Value = 2123492.8
$ws.Range("M32").Value = -2123205.8
$ws.Range("H45").Value = 460.33334
$ws.Range("I45").Value = 460.33334
$ws.Range("K45").Value = 460.33334
$ws.Range("M45").Value = -83.33334000000002
$ws.Range("H56").Value = 7000
$ws.Range("I56").Value = 7000
$ws.Range("K56").Value = 7000
$ws.Range("M56").Value = -6258
$ws.Range("H97").Value = 704.5
$ws.Range("I97").Value = 729.375
$ws.Range("K97").Value = 729.375
$ws.Range("M97").Value = -233.375
$ws.Range("H116").Value = 694.4167
$ws.Range("I116").Value = 622.2
$ws.Range("K116").Value = 622.2
$ws.Range("M116").Value = 1671.8
$ws.Range("H132").Value = 2537.1428
$ws.Range("I132").Value = 2537.1428
$ws.Range("K132").Value = 7611.428400000001
$ws.Range("M132").Value = -5081.428400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 694.4167
$ws.Range("I3").Value = 622.2
$ws.Range("K3").Value = 622.2
$ws.Range("M3").Value = -508.2
$ws.Range("H86").Value = 1876.25
$ws.Range("I86").Value = 1666
$ws.Range("J86").Value = 2507
$ws.Range("K86").Value = 1666
$ws.Range("L86").Value = 2507
$ws.Range("M86").Value = -543
$ws.Range("N86").Value = -4753
$ws.Range("H89").Value = 1876.25
$ws.Range("I89").Value = 1666
$ws.Range("J89").Value = 2507
$ws.Range("K89").Value = 8330
$ws.Range("L89").Value = 12535
$ws.Range("M89").Value = -2714
$ws.Range("N89").Value = -23767
$ws.Range("H94").Value = 400.1
$ws.Range("I94").Value = 389
$ws.Range("K94").Value = 389
$ws.Range("M94").Value = 62
$ws.Range("H105").Value = 2180.25
$ws.Range("I105").Value = 2063.1428
$ws.Range("K105").Value = 2063.1428
$ws.Range("M105").Value = -316.1428000000001
$ws.Range("H134").Value = 2066.2
$ws.Range("J134").Value = 2257
$ws.Range("L134").Value = 6771
$ws.Range("N134").Value = -11841
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 363.66666
$ws.Range("I22").Value = 363.66666
$ws.Range("K22").Value = 363.66666
$ws.Range("M22").Value = -13.66665999999998
$ws.Range("H31").Value = 2640.4285
$ws.Range("I31").Value = 2413.8333
$ws.Range("K31").Value = 2413.8333
$ws.Range("M31").Value = -2118.8333
$ws.Range("H34").Value = 2640.4285
$ws.Range("I34").Value = 2413.8333
$ws.Range("K34").Value = 2413.8333
$ws.Range("M34").Value = -2211.8333
$ws.Range("H86").Value = 9815.444
$ws.Range("I86").Value = 10294.125
$ws.Range("K86").Value = 10294.125
$ws.Range("M86").Value = -9171.125
$ws.Range("H89").Value = 9815.444
$ws.Range("I89").Value = 10294.125
$ws.Range("K89").Value = 51470.625
$ws.Range("M89").Value = -45854.625
$ws.Range("H109").Value = 56900
$ws.Range("J109").Value = 56900
$ws.Range("L109").Value = 56900
$ws.Range("N109").Value = -58980
$ws.Range("H134").Value = 3808.25
$ws.Range("I134").Value = 3033.3333
$ws.Range("J134").Value = 4583.1665
$ws.Range("K134").Value = 9099.999899999999
$ws.Range("L134").Value = 13749.4995
$ws.Range("M134").Value = -6564.999899999999
$ws.Range("N134").Value = -18819.4995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 75
$ws.Range("J34").Value = 75
$ws.Range("L34").Value = 225
$ws.Range("N34").Value = -393
$ws.Range("H55").Value = 3112
$ws.Range("I55").Value = 239
$ws.Range("J55").Value = 5985
$ws.Range("K55").Value = 717
$ws.Range("L55").Value = 17955
$ws.Range("M55").Value = -540
$ws.Range("N55").Value = -18309
$ws.Range("H113").Value = 1920.5714
$ws.Range("J113").Value = 2237.25
$ws.Range("L113").Value = 6711.75
$ws.Range("N113").Value = -11051.75
$ws.Range("H119").Value = 3332.3333
$ws.Range("I119").Value = 3332.3333
$ws.Range("K119").Value = 9996.999899999999
$ws.Range("M119").Value = -5158.999899999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 4884
$ws.Range("I10").Value = 6333.3335
$ws.Range("K10").Value = 6333.3335
$ws.Range("M10").Value = -6164.3335
$ws.Range("H15").Value = 50604.5
$ws.Range("J15").Value = 50604.5
$ws.Range("L15").Value = 50604.5
$ws.Range("N15").Value = -51180.5
$ws.Range("H81").Value = 50604.5
$ws.Range("J81").Value = 50604.5
$ws.Range("L81").Value = 50604.5
$ws.Range("N81").Value = -52600.5
$ws.Range("H84").Value = 50604.5
$ws.Range("J84").Value = 50604.5
$ws.Range("L84").Value = 151813.5
$ws.Range("N84").Value = -161797.5
$ws.Range("H113").Value = 665.6667
$ws.Range("I113").Value = 499
$ws.Range("J113").Value = 999
$ws.Range("K113").Value = 499
$ws.Range("L113").Value = 999
$ws.Range("M113").Value = 1671
$ws.Range("N113").Value = -5339
$ws.Range("H122").Value = 9999
$ws.Range("J122").Value = 9999
$ws.Range("L122").Value = 29997
$ws.Range("N122").Value = -34897
$ws.Range("H132").Value = 6757.722
$ws.Range("I132").Value = 6990.8823
$ws.Range("K132").Value = 20972.6469
$ws.Range("M132").Value = -18442.6469
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 251276
$ws.Range("J10").Value = 1701.3334
$ws.Range("L10").Value = 1701.3334
$ws.Range("N10").Value = -1981.3334
$ws.Range("H46").Value = 1780.5
$ws.Range("I46").Value = 1174.5
$ws.Range("K46").Value = 1174.5
$ws.Range("M46").Value = -986.5
$ws.Range("H55").Value = 1006.3
$ws.Range("I55").Value = 273.4
$ws.Range("J55").Value = 1739.2
$ws.Range("K55").Value = 273.4
$ws.Range("L55").Value = 1739.2
$ws.Range("M55").Value = -100.4
$ws.Range("N55").Value = -2085.2
$ws.Range("H68").Value = 2049.5
$ws.Range("I68").Value = 1600
$ws.Range("J68").Value = 2499
$ws.Range("K68").Value = 1600
$ws.Range("L68").Value = 2499
$ws.Range("M68").Value = -851
$ws.Range("N68").Value = -3997
$ws.Range("H71").Value = 2049.5
$ws.Range("I71").Value = 1600
$ws.Range("J71").Value = 2499
$ws.Range("K71").Value = 8000
$ws.Range("L71").Value = 12495
$ws.Range("M71").Value = -4256
$ws.Range("N71").Value = -19983
$ws.Range("H82").Value = 1027.8572
$ws.Range("J82").Value = 1139.2
$ws.Range("L82").Value = 1139.2
$ws.Range("N82").Value = -1861.2
$ws.Range("H85").Value = 1027.8572
$ws.Range("J85").Value = 1139.2
$ws.Range("L85").Value = 1139.2
$ws.Range("N85").Value = -3635.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 450077.5
$ws.Range("I10").Value = 155
$ws.Range("J10").Value = 900000
$ws.Range("K10").Value = 155
$ws.Range("L10").Value = 900000
$ws.Range("M10").Value = 14
$ws.Range("N10").Value = -900338
$ws.Range("H27").Value = 39989.5
$ws.Range("J27").Value = 39989.5
$ws.Range("L27").Value = 39989.5
$ws.Range("N27").Value = -40127.5
